$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the Doveton row (old row 6) into row 4
$ws.Range("A4").Value = "Doveton"
$ws.Range("B4").Value = "Holy Family Parish Doveton Catholic  100 Power Road, Doveton VIC 3177"
$ws.Range("C4").Value = "26/12/20 4:00pm-6:00pm"
$ws.Range("D4").Value = "Case attended Spanish Service"

# Move the Glen Waverley row (old row 7) into row 5
$ws.Range("A5").Value = "Glen Waverley"
$ws.Range("B5").Value = "Village Century City  285-287 Springvale Road, Glen Waverley VIC 3150"
$ws.Range("C5").Value = "28/12/20 2:45pm-5:30pm"
$ws.Range("D5").Value = "2:45pm showing of Wonder Woman 1984 (Gold Class)"

# Move the first Mordialloc row (old row 13) into row 6
$ws.Range("A6").Value = "Mordialloc"
$ws.Range("B6").Value = "Woodlands Golf Club - club bar  109 White Street Mordialloc VIC 3195"
$ws.Range("C6").Value = "23/12/20 12:30pm-1:30pm"
$ws.Range("D6").Value = "Case attended club house bar"

# Move the second Mordialloc row (old row 14) into row 7
$ws.Range("A7").Value = "Mordialloc"
$ws.Range("B7").Value = "Woodlands Golf Club - club bar  109 White Street Mordialloc VIC 3195"
$ws.Range("C7").Value = "28/12/20 4:40pm-5:15pm"
$ws.Range("D7").Value = "Case attended club house bar"

# Delete the now-duplicate trailing rows (old rows 8 through 15)
$ws.Range("A8:D15").EntireRow.Delete()
